{"js": "// tdf#112118: the trailing \"_GoBack\" bookmark (marking the last edit\n// position) moves from the final (3rd) paragraph to the first paragraph,\n// and the page-break run that used to live in the first paragraph is\n// removed. The 3rd paragraph ends up completely empty.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Build a minimal OOXML \"package\" snippet that Word can use to replace a\n// range's contents with an exact paragraph body (keeps the bookmark id\n// stable at 0, as in the original document).\nfunction packageWithParagraph(innerXml) {\n  return '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + innerXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n}\n\n// 1) First paragraph: drop the page-break run, add the _GoBack bookmark.\nfirstParagraph.getRange().insertOoxml(\n  packageWithParagraph('<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>'),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 2) Last paragraph: remove the bookmark, leaving it fully empty.\nlastParagraph.getRange().insertOoxml(\n  packageWithParagraph(''),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# tdf#112118: the trailing \"_GoBack\" bookmark (marking the last edit\n# position) moves from the final paragraph to the first paragraph, and the\n# page-break run that used to live in the first paragraph is removed. The\n# last paragraph ends up completely empty.\n\n$d = $word.ActiveDocument\n\n# 1) Remove the page-break run from the first paragraph.\n$firstRange = $d.Paragraphs(1).Range\n$firstRange.Collapse(1)        # wdCollapseStart\n$firstRange.MoveEnd(1, 1)      # wdCharacter: extend over the break char\n$firstRange.Delete()\n\n# 2) Insert the \"_GoBack\" bookmark at the (now empty) start of the first\n#    paragraph. Bookmark names are unique, so adding it here automatically\n#    removes it from wherever it used to be (the last paragraph).\n$bookmarkRange = $d.Paragraphs(1).Range\n$bookmarkRange.Collapse(1)     # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
